# Auto-generated Word COM-interop script
# Updates the date header and all 100 arithmetic problem cells
# in the single table, per the target diff.

$d = $word.ActiveDocument

# --- Update the date/weekday header paragraph ---
$d.Content.Find.Execute("2025-08-30 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-08-31 Sunday", 2) | Out-Null

# --- Update each table cell by (row, column) position ---
# Using Cell(row, col) addressing (rather than text search) because
# some old values (e.g. "65-15=") occur more than once in the table
# but map to different replacement values depending on position.
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "87-44="
$t.Cell(1, 2).Range.Text = "88-74="
$t.Cell(1, 3).Range.Text = "59+3="
$t.Cell(1, 4).Range.Text = "24+10="
$t.Cell(1, 5).Range.Text = "85-51="
$t.Cell(2, 1).Range.Text = "60-50="
$t.Cell(2, 2).Range.Text = "84-28="
$t.Cell(2, 3).Range.Text = "73-6="
$t.Cell(2, 4).Range.Text = "67-5="
$t.Cell(2, 5).Range.Text = "28+68="
$t.Cell(3, 1).Range.Text = "48-24="
$t.Cell(3, 2).Range.Text = "95+2="
$t.Cell(3, 3).Range.Text = "18+49="
$t.Cell(3, 4).Range.Text = "3+24="
$t.Cell(3, 5).Range.Text = "36-5="
$t.Cell(4, 1).Range.Text = "89+8="
$t.Cell(4, 2).Range.Text = "45-17="
$t.Cell(4, 3).Range.Text = "33-3="
$t.Cell(4, 4).Range.Text = "57-1="
$t.Cell(4, 5).Range.Text = "16+62="
$t.Cell(5, 1).Range.Text = "47-20="
$t.Cell(5, 2).Range.Text = "43-25="
$t.Cell(5, 3).Range.Text = "98-54="
$t.Cell(5, 4).Range.Text = "67-30="
$t.Cell(5, 5).Range.Text = "59-31="
$t.Cell(6, 1).Range.Text = "36+47="
$t.Cell(6, 2).Range.Text = "87-35="
$t.Cell(6, 3).Range.Text = "74-59="
$t.Cell(6, 4).Range.Text = "56-7="
$t.Cell(6, 5).Range.Text = "18+30="
$t.Cell(7, 1).Range.Text = "64+13="
$t.Cell(7, 2).Range.Text = "89-62="
$t.Cell(7, 3).Range.Text = "36-19="
$t.Cell(7, 4).Range.Text = "78-51="
$t.Cell(7, 5).Range.Text = "27+27="
$t.Cell(8, 1).Range.Text = "65-1="
$t.Cell(8, 2).Range.Text = "97-9="
$t.Cell(8, 3).Range.Text = "39-31="
$t.Cell(8, 4).Range.Text = "17+20="
$t.Cell(8, 5).Range.Text = "22-9="
$t.Cell(9, 1).Range.Text = "96-31="
$t.Cell(9, 2).Range.Text = "33-3="
$t.Cell(9, 3).Range.Text = "71-11="
$t.Cell(9, 4).Range.Text = "61+34="
$t.Cell(9, 5).Range.Text = "33-6="
$t.Cell(10, 1).Range.Text = "93-77="
$t.Cell(10, 2).Range.Text = "67-1="
$t.Cell(10, 3).Range.Text = "51-8="
$t.Cell(10, 4).Range.Text = "84-9="
$t.Cell(10, 5).Range.Text = "66-13="
$t.Cell(11, 1).Range.Text = "91-82="
$t.Cell(11, 2).Range.Text = "3+64="
$t.Cell(11, 3).Range.Text = "91-89="
$t.Cell(11, 4).Range.Text = "72+19="
$t.Cell(11, 5).Range.Text = "21+74="
$t.Cell(12, 1).Range.Text = "56-20="
$t.Cell(12, 2).Range.Text = "57-9="
$t.Cell(12, 3).Range.Text = "74-58="
$t.Cell(12, 4).Range.Text = "79+19="
$t.Cell(12, 5).Range.Text = "85-82="
$t.Cell(13, 1).Range.Text = "11+46="
$t.Cell(13, 2).Range.Text = "33+59="
$t.Cell(13, 3).Range.Text = "23+60="
$t.Cell(13, 4).Range.Text = "92-14="
$t.Cell(13, 5).Range.Text = "76+19="
$t.Cell(14, 1).Range.Text = "30+7="
$t.Cell(14, 2).Range.Text = "32-11="
$t.Cell(14, 3).Range.Text = "25-17="
$t.Cell(14, 4).Range.Text = "92-78="
$t.Cell(14, 5).Range.Text = "41-2="
$t.Cell(15, 1).Range.Text = "88-83="
$t.Cell(15, 2).Range.Text = "59-35="
$t.Cell(15, 3).Range.Text = "22+41="
$t.Cell(15, 4).Range.Text = "90-36="
$t.Cell(15, 5).Range.Text = "21+52="
$t.Cell(16, 1).Range.Text = "68+6="
$t.Cell(16, 2).Range.Text = "55+27="
$t.Cell(16, 3).Range.Text = "95-66="
$t.Cell(16, 4).Range.Text = "81-1="
$t.Cell(16, 5).Range.Text = "37+0="
$t.Cell(17, 1).Range.Text = "35+23="
$t.Cell(17, 2).Range.Text = "3+74="
$t.Cell(17, 3).Range.Text = "74+9="
$t.Cell(17, 4).Range.Text = "92-27="
$t.Cell(17, 5).Range.Text = "20-13="
$t.Cell(18, 1).Range.Text = "25+62="
$t.Cell(18, 2).Range.Text = "30+44="
$t.Cell(18, 3).Range.Text = "11+69="
$t.Cell(18, 4).Range.Text = "7+18="
$t.Cell(18, 5).Range.Text = "37+58="
$t.Cell(19, 1).Range.Text = "63-37="
$t.Cell(19, 2).Range.Text = "79-27="
$t.Cell(19, 3).Range.Text = "77+3="
$t.Cell(19, 4).Range.Text = "20+12="
$t.Cell(19, 5).Range.Text = "44-37="
$t.Cell(20, 1).Range.Text = "6+44="
$t.Cell(20, 2).Range.Text = "70-16="
$t.Cell(20, 3).Range.Text = "82-31="
$t.Cell(20, 4).Range.Text = "10+30="
$t.Cell(20, 5).Range.Text = "78-52="

Write-Output "Done"
